{"js": "// Append the missing SRS section numbers to the two journal entries that\n// only say \"SRS Section\" / \"SRS Section \" with nothing after it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet firstDone = false;\nlet secondDone = false;\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n\n  // First occurrence: \"SRS Section \" (trailing space, nothing after it yet).\n  if (!firstDone && text === \"SRS Section \") {\n    paragraph.insertText(\"1.1, \", Word.InsertLocation.end);\n    paragraph.insertText(\"1.3, 2.2\", Word.InsertLocation.end);\n    firstDone = true;\n    continue;\n  }\n\n  // Second occurrence: \"SRS Section\" (no trailing space, nothing after it yet).\n  if (!secondDone && text === \"SRS Section\") {\n    paragraph.insertText(\" 4.3, 4.4\", Word.InsertLocation.end);\n    secondDone = true;\n    continue;\n  }\n\n  if (firstDone && secondDone) break;\n}\n\nawait context.sync();\n", "ps1": "# Append the missing SRS section numbers to the two journal entries that\n# only say \"SRS Section\" / \"SRS Section \" with nothing after it.\n$d = $word.ActiveDocument\n$cr = [char]13\n\n$firstDone = $false\n$secondDone = $false\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    if (-not $firstDone -and $t -eq (\"SRS Section \" + $cr)) {\n        $p.Range.InsertAfter(\"1.1, \")\n        $p.Range.InsertAfter(\"1.3, 2.2\")\n        $firstDone = $true\n        continue\n    }\n\n    if (-not $secondDone -and $t -eq (\"SRS Section\" + $cr)) {\n        $p.Range.InsertAfter(\" 4.3, 4.4\")\n        $secondDone = $true\n        continue\n    }\n}\n"}
